# Apply cryptos list update (Mon May 27 18:52:59 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.944.99"
$ws.Range("E2").Value = "  +1.63%  "

# Row 3
$ws.Range("D3").Value = "3.926.06"
$ws.Range("E3").Value = "  +2.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.97%  "

# Row 7
$ws.Range("D7").Value = "3.927.82"
$ws.Range("E7").Value = "  +2.09%  "

# Row 8
$ws.Range("E8").Value = "  +0.22%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.13%  "

# Row 10
$ws.Range("E10").Value = "  +1.08%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.78%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.469"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.23%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.06%  "

# Row 15
$ws.Range("D15").Value = "4.590.45"
$ws.Range("E15").Value = "  +2.25%  "

# Row 16
$ws.Range("D16").Value = "3.933.12"
$ws.Range("E16").Value = "  +1.72%  "

# Row 17
$ws.Range("D17").Value = "69.967.70"
$ws.Range("E17").Value = "  +1.45%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.33%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.47%  "

# Row 20
$ws.Range("E20").Value = "  -0.84%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.67%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.96%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.749"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.23%  "

# Row 24
$ws.Range("E24").Value = "  +5.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.47%  "

# Row 26
$ws.Range("E26").Value = "  +2.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.95%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.60%  "

# Row 29
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.62%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.09%  "

# Row 32
$ws.Range("D32").Value = "4.079.30"
$ws.Range("E32").Value = "  +2.18%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.91%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.13%  "

# Row 35
$ws.Range("D35").Value = "3.891.89"
$ws.Range("E35").Value = "  +2.62%  "

# Row 36
$ws.Range("E36").Value = "  +1.35%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.22%  "

# Row 38
$ws.Range("E38").Value = "  +2.20%  "

# Row 39
$ws.Range("E39").Value = "  +1.76%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.00%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.18%  "

# Row 42
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.330"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.41%  "

# Row 43
$ws.Range("E43").Value = "  +8.17%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "437.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.18%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.42%  "

# Row 47
$ws.Range("E47").Value = "  +0.03%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0371"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.30%  "

# Row 49
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000277"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +23.20%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.18%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.24%  "
